$d = $word.ActiveDocument

# Replace paragraph 2's content first (it is not yet the last paragraph, so
# InsertXML cleanly swaps its runs/oMath without leaving a stray empty paragraph).
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:r><w:t>2</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t>若</w:t></w:r><m:oMath xmlns:mml="http://www.w3.org/1998/Math/MathML"><m:bar><m:barPr><m:pos m:val="top"/></m:barPr><m:e><m:r><m:t>BD</m:t></m:r></m:e></m:bar><m:r><m:t>：</m:t></m:r><m:bar><m:barPr><m:pos m:val="top"/></m:barPr><m:e><m:r><m:t>CD</m:t></m:r></m:e></m:bar></m:oMath><w:r><w:t>＝</w:t></w:r><w:r><w:t>2</w:t></w:r><w:r><w:t>：</w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t>，</w:t></w:r><m:oMath xmlns:mml="http://www.w3.org/1998/Math/MathML"><m:bar><m:barPr><m:pos m:val="top"/></m:barPr><m:e><m:r><m:t>AE</m:t></m:r></m:e></m:bar><m:r><m:t>：</m:t></m:r><m:bar><m:barPr><m:pos m:val="top"/></m:barPr><m:e><m:r><m:t>DE</m:t></m:r></m:e></m:bar></m:oMath><w:r><w:t>＝</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>：</w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t>，</w:t></w:r><w:r><w:t>則</w:t></w:r><w:r><w:t>△</w:t></w:r><w:r><w:t>A</w:t></w:r><w:r><w:t>B</w:t></w:r><w:r><w:t>E</w:t></w:r><w:r><w:t>面</w:t></w:r><w:r><w:t>積</w:t></w:r><w:r><w:t>：</w:t></w:r><w:r><w:t>△</w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:t>D</w:t></w:r><w:r><w:t>E</w:t></w:r><w:r><w:t>面</w:t></w:r><w:r><w:t>積</w:t></w:r><w:r><w:t>＝</w:t></w:r><w:r><w:t>【</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>：</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t>0</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>】</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

# Drop the old paragraphs 3 ("測試線段題目…") and 4 ("測試跟號…"); their bar()
# sub-expressions were already folded into the rebuilt paragraph 2 above.
$p3 = $d.Paragraphs.Item(3)
$p4 = $d.Paragraphs.Item(4)
$d.Range($p3.Range.Start, $p4.Range.End).Delete()

# Replace paragraph 1's content last.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML('<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><w:r><w:t>1</w:t></w:r><w:r><w:t>.</w:t></w:r><m:oMath xmlns:mml="http://www.w3.org/1998/Math/MathML"><m:f><m:fPr><m:type m:val="bar"/></m:fPr><m:num><m:r><m:t>A</m:t></m:r></m:num><m:den><m:r><m:t>Z</m:t></m:r></m:den></m:f></m:oMath><w:r><w:t>C</w:t></w:r><w:r><w:t>L</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>請</w:t></w:r><w:r><w:t>問</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>b</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>=</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
